$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set cell values in the exact order needed so the shared-string table ---
# --- is built in the same sequence as the authored workbook.               ---

# Header row (row 1) - first occurrences of header labels
$ws.Cells.Item(1, 3).Value = "Type"
$ws.Cells.Item(1, 5).Value = "Case"
$ws.Cells.Item(1, 2).Value = "const/let"
$ws.Cells.Item(1, 4).Value = "Location"

# Row 2 (playerMoney) - first occurrences of the repeated value set
$ws.Cells.Item(2, 2).Value = "let"
$ws.Cells.Item(2, 3).Value = "number"
$ws.Cells.Item(2, 4).Value = "player"
$ws.Cells.Item(2, 5).Value = "camel"

# Remaining header cells
$ws.Cells.Item(1, 1).Value = "Current Name"
$ws.Cells.Item(1, 6).Value = "New Name"

# Column A values (variable names) for the rest of the rows
$ws.Cells.Item(2, 1).Value = "playerMoney"
$ws.Cells.Item(3, 1).Value = "lastPlayerMoney"
$ws.Cells.Item(4, 1).Value = "playerRoll"
$ws.Cells.Item(5, 1).Value = "playerPoint"

# Usage column (G) header + notes
$ws.Cells.Item(1, 7).Value = "Usage"
$ws.Cells.Item(2, 7).Value = "Track player's money"
$ws.Cells.Item(4, 7).Value = "Track last dice roll"
$ws.Cells.Item(5, 7).Value = "Track current point"

# Fill remaining repeated const/let, Type, Location, Case cells for rows 3-5
for ($r = 3; $r -le 5; $r++) {
    $ws.Cells.Item($r, 2).Value = "let"
    $ws.Cells.Item($r, 3).Value = "number"
    $ws.Cells.Item($r, 4).Value = "player"
    $ws.Cells.Item($r, 5).Value = "camel"
}

# --- Column widths (approximate best-fit widths from the authored file) ---
$ws.Columns.Item(1).ColumnWidth = 14.5
$ws.Columns.Item(2).ColumnWidth = 10.333333333333332
$ws.Columns.Item(3).ColumnWidth = 9.5
$ws.Columns.Item(6).ColumnWidth = 11.833333333333332
$ws.Columns.Item(7).ColumnWidth = 17

# --- Turn the range into an Excel Table ("Table1") ---
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:G14"), $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = "TableStyleLight8"

# --- Selection shown when the workbook was saved ---
$ws.Range("G6").Select() | Out-Null
